$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the misspelled "Thrid" -> "Third"
$ws.Range("A2").Value = "Third"

# Update TimeCol values (column D) for rows 3, 4, 5
$ws.Range("D3").Value = 0.01707175925925926
$ws.Range("D4").Value = 0.5163888888888889
$ws.Range("D5").Value = 0.47890046296296296

# Match the saved selection state (active cell D3)
$ws.Range("D3").Select()

$wb.Save()
